$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells that will receive numeric-looking strings
# so they remain text (matching the inline-string semantics of the source).
$textCells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","E24","E25","D27","E27","D39","E39","D40","E40","D41","E41","E42","E43","D44","E44","D45","E45","D46","E46","D48","E48","E49","D50","E50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply cell value changes from the diff
$ws.Range("D2").Value = '305.98'
$ws.Range("E2").Value = '0.83%'
$ws.Range("D3").Value = '36.14'
$ws.Range("E3").Value = '-2.72%'
$ws.Range("E4").Value = '2.43%'
$ws.Range("D5").Value = '0.07872'
$ws.Range("E5").Value = '0.50%'
$ws.Range("D6").Value = '2.140'
$ws.Range("E6").Value = '-3.67%'
$ws.Range("D7").Value = '7.942'
$ws.Range("E7").Value = '-1.05%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '0.9179'
$ws.Range("E8").Value = '0.33%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '0.09650'
$ws.Range("E9").Value = '-0.68%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1866'
$ws.Range("E10").Value = '-1.09%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.08727'
$ws.Range("E11").Value = '1.73%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03580'
$ws.Range("E12").Value = '0.97%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09930'
$ws.Range("E13").Value = '-0.25%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001428'
$ws.Range("E14").Value = '-4.65%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.005621'
$ws.Range("E15").Value = '-0.67%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.456'
$ws.Range("E16").Value = '-0.15%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '4.110'
$ws.Range("E17").Value = '1.74%'
$ws.Range("D18").Value = '2.712'
$ws.Range("E18").Value = '19.85%'
$ws.Range("D19").Value = '0.3395'
$ws.Range("E19").Value = '-1.94%'
$ws.Range("D20").Value = '0.1328'
$ws.Range("E20").Value = '2.14%'
$ws.Range("D21").Value = '5.175'
$ws.Range("E21").Value = '8.17%'
$ws.Range("D22").Value = '0.2020'
$ws.Range("E22").Value = '-12.03%'
$ws.Range("D23").Value = '0.04556'
$ws.Range("E23").Value = '-1.04%'
$ws.Range("E24").Value = '5.37%'
$ws.Range("E25").Value = '0.33%'
$ws.Range("D27").Value = '0.0004749'
$ws.Range("E27").Value = '-0.04%'
$ws.Range("D39").Value = '0.01851'
$ws.Range("E39").Value = '3.72%'
$ws.Range("D40").Value = '0.04761'
$ws.Range("E40").Value = '0.21%'
$ws.Range("D41").Value = '0.007484'
$ws.Range("E41").Value = '-6.63%'
$ws.Range("E42").Value = '0.47%'
$ws.Range("E43").Value = '0.68%'
$ws.Range("D44").Value = '0.002249'
$ws.Range("E44").Value = '4.05%'
$ws.Range("D45").Value = '0.01104'
$ws.Range("E45").Value = '14.80%'
$ws.Range("D46").Value = '0.00006319'
$ws.Range("E46").Value = '3.24%'
$ws.Range("D48").Value = '0.0005800'
$ws.Range("E48").Value = '-0.01%'
$ws.Range("E49").Value = '524.66%'
$ws.Range("D50").Value = '0.002000'
$ws.Range("E50").Value = '-25.67%'
